$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_5_0_0"
$ws.Range("B2").Value = 0.2429408329254606
$ws.Range("C2").Value = 0.1514423548127821
$ws.Range("D2").Value = 0.3481670994364133
$ws.Range("E2").Value = 0.212135861104836
$ws.Range("F2").Value = 0.8378413915634155
$ws.Range("G2").Value = 1.449533343315125
$ws.Range("H2").Value = 0.4810515940189362
$ws.Range("I2").Value = 0.9937768578529358

$ws.Range("A3").Value = "model_5_0_1"
$ws.Range("B3").Value = 0.2512561413512001
$ws.Range("C3").Value = 0.1305860278065826
$ws.Range("D3").Value = -0.2146526476029917
$ws.Range("E3").Value = 0.04221911415414981
$ws.Range("F3").Value = 0.8286388516426086
$ws.Range("G3").Value = 1.485160946846008
$ws.Range("H3").Value = 0.896411657333374
$ws.Range("I3").Value = 1.208102345466614

$ws.Range("A4").Value = "model_5_0_2"
$ws.Range("B4").Value = 0.4215613069722006
$ws.Range("C4").Value = 0.2176362758038078
$ws.Range("D4").Value = -0.208595957848797
$ws.Range("E4").Value = 0.1062995981093575
$ws.Range("F4").Value = 0.6401611566543579
$ws.Range("G4").Value = 1.336458802223206
$ws.Range("H4").Value = 0.8919417858123779
$ws.Range("I4").Value = 1.127274036407471

$ws.Range("A5").Value = "model_5_0_10"
$ws.Range("B5").Value = 0.6547333304004781
$ws.Range("C5").Value = 0.4678236617498618
$ws.Range("D5").Value = -0.8205131655388846
$ws.Range("E5").Value = 0.1171956668846853
$ws.Range("F5").Value = 0.3821084499359131
$ws.Range("G5").Value = 0.9090807437896729
$ws.Range("H5").Value = 1.343535661697388
$ws.Range("I5").Value = 1.113530158996582

$ws.Range("A6").Value = "model_5_0_11"
$ws.Range("B6").Value = 0.6547881192562768
$ws.Range("C6").Value = 0.4711413548303569
$ws.Range("D6").Value = -0.8447227508122814
$ws.Range("E6").Value = 0.1129081845513944
$ws.Range("F6").Value = 0.3820478320121765
$ws.Range("G6").Value = 0.9034133553504944
$ws.Range("H6").Value = 1.36140239238739
$ws.Range("I6").Value = 1.118938446044922

$ws.Range("A7").Value = "model_5_0_12"
$ws.Range("B7").Value = 0.6554996122614334
$ws.Range("C7").Value = 0.4732600807310715
$ws.Range("D7").Value = -0.8577507561234246
$ws.Range("E7").Value = 0.1108409886899749
$ws.Range("F7").Value = 0.3812603950500488
$ws.Range("G7").Value = 0.8997940421104431
$ws.Range("H7").Value = 1.371016979217529
$ws.Range("I7").Value = 1.121545791625977

$ws.Range("A8").Value = "model_5_0_8"
$ws.Range("B8").Value = 0.6556823423532738
$ws.Range("C8").Value = 0.4637005135354108
$ws.Range("D8").Value = -0.766073887651445
$ws.Range("E8").Value = 0.1292283746967942
$ws.Range("F8").Value = 0.3810581862926483
$ws.Range("G8").Value = 0.9161240458488464
$ws.Range("H8").Value = 1.303359508514404
$ws.Range("I8").Value = 1.098352789878845

$ws.Range("A9").Value = "model_5_0_13"
$ws.Range("B9").Value = 0.6571722164373955
$ws.Range("C9").Value = 0.4761530846898955
$ws.Range("D9").Value = -0.862351234753975
$ws.Range("E9").Value = 0.1116482564195721
$ws.Range("F9").Value = 0.3794093132019043
$ws.Range("G9").Value = 0.8948521018028259
$ws.Range("H9").Value = 1.374412178993225
$ws.Range("I9").Value = 1.120527505874634

$ws.Range("A10").Value = "model_5_0_7"
$ws.Range("B10").Value = 0.6572333450553101
$ws.Range("C10").Value = 0.4406539254847771
$ws.Range("D10").Value = -0.6073738523754888
$ws.Range("E10").Value = 0.156400858159532
$ws.Range("F10").Value = 0.3793416917324066
$ws.Range("G10").Value = 0.9554928541183472
$ws.Range("H10").Value = 1.186239242553711
$ws.Range("I10").Value = 1.064078450202942

$ws.Range("A11").Value = "model_5_0_24"
$ws.Range("B11").Value = 0.6572425462059982
$ws.Range("C11").Value = 0.4920657159085137
$ws.Range("D11").Value = -1.036689706991371
$ws.Range("E11").Value = 0.07505748579889249
$ws.Range("F11").Value = 0.37933149933815
$ws.Range("G11").Value = 0.8676695823669434
$ws.Range("H11").Value = 1.503073573112488
$ws.Range("I11").Value = 1.166681408882141

$ws.Range("A12").Value = "model_5_0_9"
$ws.Range("B12").Value = 0.6580715254567583
$ws.Range("C12").Value = 0.4673827882810914
$ws.Range("D12").Value = -0.7672782302719907
$ws.Range("E12").Value = 0.1315377741221033
$ws.Range("F12").Value = 0.3784140944480896
$ws.Range("G12").Value = 0.9098337888717651
$ws.Range("H12").Value = 1.304248452186584
$ws.Range("I12").Value = 1.095439791679382

$ws.Range("A13").Value = "model_5_0_16"
$ws.Range("B13").Value = 0.6582109063610959
$ws.Range("C13").Value = 0.4816145732107079
$ws.Range("D13").Value = -0.9081562989287406
$ws.Range("E13").Value = 0.102951653418262
$ws.Range("F13").Value = 0.3782598376274109
$ws.Range("G13").Value = 0.8855226039886475
$ws.Range("H13").Value = 1.408216238021851
$ws.Range("I13").Value = 1.131497025489807

$ws.Range("A14").Value = "model_5_0_17"
$ws.Range("B14").Value = 0.658344975524519
$ws.Range("C14").Value = 0.4828270785343209
$ws.Range("D14").Value = -0.9229057309604229
$ws.Range("E14").Value = 0.09976061174113304
$ws.Range("F14").Value = 0.3781114220619202
$ws.Range("G14").Value = 0.8834514021873474
$ws.Range("H14").Value = 1.419101238250732
$ws.Range("I14").Value = 1.13552188873291

$ws.Range("A15").Value = "model_5_0_14"
$ws.Range("B15").Value = 0.6585396417095928
$ws.Range("C15").Value = 0.4787853935568881
$ws.Range("D15").Value = -0.8697041448930178
$ws.Range("E15").Value = 0.1115106122315868
$ws.Range("F15").Value = 0.3778960108757019
$ws.Range("G15").Value = 0.8903555274009705
$ws.Range("H15").Value = 1.379838585853577
$ws.Range("I15").Value = 1.120701193809509

$ws.Range("A16").Value = "model_5_0_15"
$ws.Range("B16").Value = 0.6587945774119868
$ws.Range("C16").Value = 0.4801459548104384
$ws.Range("D16").Value = -0.885013989863187
$ws.Range("E16").Value = 0.1082717353455406
$ws.Range("F16").Value = 0.3776138424873352
$ws.Range("G16").Value = 0.8880313038825989
$ws.Range("H16").Value = 1.3911372423172
$ws.Range("I16").Value = 1.124786376953125

$ws.Range("A17").Value = "model_5_0_23"
$ws.Range("B17").Value = 0.6588171345639016
$ws.Range("C17").Value = 0.4917423763174339
$ws.Range("D17").Value = -1.009279900591719
$ws.Range("E17").Value = 0.08237193141438504
$ws.Range("F17").Value = 0.3775889277458191
$ws.Range("G17").Value = 0.8682219982147217
$ws.Range("H17").Value = 1.482845306396484
$ws.Range("I17").Value = 1.157455325126648

$ws.Range("A18").Value = "model_5_0_22"
$ws.Range("B18").Value = 0.6588686562906014
$ws.Range("C18").Value = 0.4903799895028725
$ws.Range("D18").Value = -0.9953053735012991
$ws.Range("E18").Value = 0.08524180993574759
$ws.Range("F18").Value = 0.3775318562984467
$ws.Range("G18").Value = 0.8705492615699768
$ws.Range("H18").Value = 1.472532153129578
$ws.Range("I18").Value = 1.153835415840149

$ws.Range("A19").Value = "model_5_0_21"
$ws.Range("B19").Value = 0.6592121508394024
$ws.Range("C19").Value = 0.4894065237535035
$ws.Range("D19").Value = -0.9780189514291955
$ws.Range("E19").Value = 0.08930335298107595
$ws.Range("F19").Value = 0.3771517276763916
$ws.Range("G19").Value = 0.8722121715545654
$ws.Range("H19").Value = 1.459774732589722
$ws.Range("I19").Value = 1.148712515830994

$ws.Range("A20").Value = "model_5_0_20"
$ws.Range("B20").Value = 0.6592506943321383
$ws.Range("C20").Value = 0.4880430000219972
$ws.Range("D20").Value = -0.9629455070698305
$ws.Range("E20").Value = 0.09247642242558995
$ws.Range("F20").Value = 0.3771090507507324
$ws.Range("G20").Value = 0.8745413422584534
$ws.Range("H20").Value = 1.448650479316711
$ws.Range("I20").Value = 1.144709825515747

$ws.Range("A21").Value = "model_5_0_19"
$ws.Range("B21").Value = 0.659483002774766
$ws.Range("C21").Value = 0.4866504942373159
$ws.Range("D21").Value = -0.9449406253660155
$ws.Range("E21").Value = 0.09643503239882434
$ws.Range("F21").Value = 0.3768519759178162
$ws.Range("G21").Value = 0.8769200444221497
$ws.Range("H21").Value = 1.435362935066223
$ws.Range("I21").Value = 1.139716863632202

$ws.Range("A22").Value = "model_5_0_18"
$ws.Range("B22").Value = 0.6597552101276527
$ws.Range("C22").Value = 0.4850706010602949
$ws.Range("D22").Value = -0.926270823051853
$ws.Range("E22").Value = 0.1004430343495436
$ws.Range("F22").Value = 0.376550704240799
$ws.Range("G22").Value = 0.8796189427375793
$ws.Range("H22").Value = 1.421584725379944
$ws.Range("I22").Value = 1.134661316871643

$ws.Range("A23").Value = "model_5_0_3"
$ws.Range("B23").Value = 0.6710990045298177
$ws.Range("C23").Value = 0.464259267051516
$ws.Range("D23").Value = 0.141150367226888
$ws.Range("E23").Value = 0.3794182273120696
$ws.Range("F23").Value = 0.3639964759349823
$ws.Range("G23").Value = 0.9151694774627686
$ws.Range("H23").Value = 0.6338295936584473
$ws.Range("I23").Value = 0.7827743887901306

$ws.Range("A24").Value = "model_5_0_6"
$ws.Range("B24").Value = 0.679616229954948
$ws.Range("C24").Value = 0.4228327936253783
$ws.Range("D24").Value = -0.2253337214069369
$ws.Range("E24").Value = 0.2488118830841948
$ws.Range("F24").Value = 0.3545704185962677
$ws.Range("G24").Value = 0.9859356284141541
$ws.Range("H24").Value = 0.9042942523956299
$ws.Range("I24").Value = 0.9475153684616089

$ws.Range("A25").Value = "model_5_0_4"
$ws.Range("B25").Value = 0.697270525740101
$ws.Range("C25").Value = 0.4804769280168062
$ws.Range("D25").Value = 0.2325617257503795
$ws.Range("E25").Value = 0.4162147807610748
$ws.Range("F25").Value = 0.3350323140621185
$ws.Range("G25").Value = 0.8874659538269043
$ws.Range("H25").Value = 0.5663681626319885
$ws.Range("I25").Value = 0.7363608479499817

$ws.Range("A26").Value = "model_5_0_5"
$ws.Range("B26").Value = 0.7195745727398117
$ws.Range("C26").Value = 0.4690497611302625
$ws.Range("D26").Value = 0.1711086737449582
$ws.Range("E26").Value = 0.3911017596024546
$ws.Range("F26").Value = 0.3103483021259308
$ws.Range("G26").Value = 0.9069862365722656
$ws.Range("H26").Value = 0.6117204427719116
$ws.Range("I26").Value = 0.7680372595787048

# Ensure newly added rows (25-26) in column A inherit the existing label style
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A25:A26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Write-Host "Applied updates"
